$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 743
$ws.Range("D3").Value = 120
$ws.Range("D4").Value = 340
$ws.Range("D5").Value = 150
$ws.Range("D6").Value = 157
$ws.Range("D7").Value = 1620
$ws.Range("D8").Value = 650
$ws.Range("D9").Value = 145
